$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text before writing, so numeric-looking
# strings like "1.00" or "0.0000188" are not auto-converted to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "74.800.19"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "2.845.79"
$ws.Range("E3").Value = "  +9.77%  "

$ws.Range("D5").Value = "188.65"
$ws.Range("E5").Value = "  +1.76%  "

$ws.Range("D6").Value = "600.34"
$ws.Range("E6").Value = "  +3.40%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "0.559"
$ws.Range("E8").Value = "  +4.76%  "

$ws.Range("E9").Value = "  -6.44%  "

$ws.Range("D10").Value = "2.842.99"
$ws.Range("E10").Value = "  +9.68%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").Value = "0.371"
$ws.Range("E12").Value = "  +3.26%  "

$ws.Range("D13").Value = "4.91"
$ws.Range("E13").Value = "  +2.72%  "

$ws.Range("D14").Value = "3.369.31"
$ws.Range("E14").Value = "  +10.58%  "

$ws.Range("D15").Value = "75.093.47"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").Value = "27.19"
$ws.Range("E16").Value = "  +3.74%  "

$ws.Range("D17").Value = "0.0000188"
$ws.Range("E17").Value = "  -2.19%  "

$ws.Range("D18").Value = "2.840.48"
$ws.Range("E18").Value = "  +9.46%  "

$ws.Range("D19").Value = "9.10"
$ws.Range("E19").Value = "  +2.58%  "

$ws.Range("D20").Value = "12.44"
$ws.Range("E20").Value = "  +5.65%  "

$ws.Range("D21").Value = "376.47"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").Value = "2.27"
$ws.Range("E22").Value = "  -2.27%  "

$ws.Range("E23").Value = "  +2.03%  "

$ws.Range("D24").Value = "6.15"
$ws.Range("E24").Value = "  -1.63%  "

$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").Value = "70.81"
$ws.Range("E26").Value = "  +1.41%  "

$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").Value = "4.21"
$ws.Range("E27").Value = "  +0.94%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.983.50"
$ws.Range("E28").Value = "  +9.82%  "

$ws.Range("D29").Value = "9.65"
$ws.Range("E29").Value = "  +4.26%  "

$ws.Range("E30").Value = "  +10.00%  "

$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("D32").Value = "531.41"
$ws.Range("E32").Value = "  +3.90%  "

$ws.Range("E33").Value = "  +2.92%  "

$ws.Range("D34").Value = "7.97"
$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").Value = "1.82"
$ws.Range("E35").Value = "  +5.25%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "20.26"
$ws.Range("E37").Value = "  +5.78%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.121"
$ws.Range("E38").Value = "  +1.74%  "

$ws.Range("D39").Value = "162.08"
$ws.Range("E39").Value = "  +1.57%  "

$ws.Range("E40").Value = "  -0.56%  "

$ws.Range("D41").Value = "184.94"
$ws.Range("E41").Value = "  +21.66%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").Value = "5.09"
$ws.Range("E43").Value = "  +3.18%  "

$ws.Range("D44").Value = "0.342"
$ws.Range("E44").Value = "  +6.19%  "

$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("E46").Value = "  +5.82%  "

$ws.Range("D47").Value = "39.66"
$ws.Range("E47").Value = "  +1.98%  "

$ws.Range("E48").Value = "  -3.74%  "

$ws.Range("D49").Value = "0.0864"
$ws.Range("E49").Value = "  +5.21%  "

$ws.Range("D50").Value = "0.574"
$ws.Range("E50").Value = "  +10.03%  "

$ws.Range("D51").Value = "3.77"
$ws.Range("E51").Value = "  +3.82%  "

# Restore the default style on the Price column so no stray
# cell-format attributes are left behind (values remain text).
$priceRange.Style = "Normal"
